# Auto-generated script: adds 2022-12-10 daily crime counts
# across Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 6905
$ws.Range("I3").Value = 7169
$ws.Range("I4").Value = 1646
$ws.Range("I5").Value = 672
$ws.Range("I6").Value = 8418
$ws.Range("I7").Value = 24810

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 196
$ws.Range("I6").Value = 179
$ws.Range("I7").Value = 776
$ws.Range("I8").Value = 1483
$ws.Range("I9").Value = 129
$ws.Range("I11").Value = 378
$ws.Range("I12").Value = 63
$ws.Range("I15").Value = 288
$ws.Range("I17").Value = 36
$ws.Range("I19").Value = 698
$ws.Range("I20").Value = 612
$ws.Range("I21").Value = 111
$ws.Range("I23").Value = 244
$ws.Range("I24").Value = 68
$ws.Range("I25").Value = 130
$ws.Range("I27").Value = 218
$ws.Range("I29").Value = 1486
$ws.Range("I31").Value = 248
$ws.Range("I32").Value = 29
$ws.Range("I33").Value = 1094
$ws.Range("I37").Value = 766
$ws.Range("I42").Value = 932
$ws.Range("I45").Value = 49
$ws.Range("I48").Value = 313
$ws.Range("I49").Value = 164
$ws.Range("I51").Value = 291
$ws.Range("I52").Value = 560
$ws.Range("I53").Value = 282
$ws.Range("I54").Value = 488
$ws.Range("I55").Value = 286
$ws.Range("I57").Value = 100
$ws.Range("I64").Value = 198
$ws.Range("I66").Value = 72
$ws.Range("I67").Value = 940
$ws.Range("I73").Value = 223
$ws.Range("I76").Value = 352
$ws.Range("I79").Value = 711
$ws.Range("I80").Value = 78
$ws.Range("I83").Value = 533
$ws.Range("I84").Value = 219
$ws.Range("I85").Value = 1109
$ws.Range("I87").Value = 63
$ws.Range("I89").Value = 293
$ws.Range("I90").Value = 321
$ws.Range("I91").Value = 260
$ws.Range("I96").Value = 288
$ws.Range("I97").Value = 222
$ws.Range("I98").Value = 182
$ws.Range("I99").Value = 434
$ws.Range("I101").Value = 24810

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 317
$ws.Range("I7").Value = 1109

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 182
$ws.Range("I7").Value = 560

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I4").Value = 38
$ws.Range("I7").Value = 378

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 428
$ws.Range("I6").Value = 479
$ws.Range("I7").Value = 1483

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 138
$ws.Range("I7").Value = 282

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I6").Value = 209
$ws.Range("I7").Value = 776

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 293

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I3").Value = 65
$ws.Range("I6").Value = 116
$ws.Range("I7").Value = 288

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I6").Value = 233
$ws.Range("I7").Value = 766

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I6").Value = 112
$ws.Range("I7").Value = 434

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 350
$ws.Range("I6").Value = 282
$ws.Range("I7").Value = 940

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 75
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 248

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 79
$ws.Range("I7").Value = 219

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 177
$ws.Range("I7").Value = 533

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 351
$ws.Range("I7").Value = 1094

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I4").Value = 20
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I2").Value = 105
$ws.Range("I7").Value = 488

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 439
$ws.Range("I3").Value = 510
$ws.Range("I4").Value = 79
$ws.Range("I7").Value = 1486

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I6").Value = 224
$ws.Range("I7").Value = 698

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 59
$ws.Range("I6").Value = 162
$ws.Range("I7").Value = 313

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 160
$ws.Range("I7").Value = 352

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I6").Value = 370
$ws.Range("I7").Value = 932

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 86
$ws.Range("I3").Value = 89
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 87
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("I2").Value = 82
$ws.Range("I7").Value = 260

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 111

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 208
$ws.Range("I3").Value = 231
$ws.Range("I5").Value = 28
$ws.Range("I7").Value = 711

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 173
$ws.Range("I6").Value = 213
$ws.Range("I7").Value = 612

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 49
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 68
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 288

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I6").Value = 118
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 73
$ws.Range("I7").Value = 223

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 74
$ws.Range("I3").Value = 61
$ws.Range("I7").Value = 196

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("I3").Value = 4
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 61
$ws.Range("I3").Value = 42
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 218

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 103
$ws.Range("I3").Value = 83
$ws.Range("I7").Value = 321

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 291

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I3").Value = 24
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 49

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I6").Value = 37
$ws.Range("I7").Value = 63
